# Updated legacy GSC export data:
# The "Chart" sheet's row 2 (the placeholder "2025-10-01" row with no real
# data yet) is removed; every subsequent row shifts up by one, and the
# sheet now reports one extra day's worth of data that was previously
# missing a row (2025-10-02 .. 2025-12-29), ending at row 90 instead of 91.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Delete entire row 2 and shift the remaining rows up.
$ws.Rows.Item(2).Delete()
